# Remove the "Version 1.0 vom "2020-10-20"" paragraph (pStyle "Date")
# that used to follow the "... Hochschulen" subtitle, per the commit's
# removal of the date/version line from the title page.
$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Date") {
        $p.Range.Delete()
        break
    }
}
